{"js": "// Replace the date line and all the two-digit-by-two-digit multiplication\n// prompts in the table with the new values from the commit.\nconst replacements = [\n  [\"2025-05-14 Wednesday\", \"2025-05-15 Thursday\"],\n  [\"17\u00d780=\", \"99\u00d799=\"],\n  [\"35\u00d762=\", \"84\u00d717=\"],\n  [\"74\u00d772=\", \"39\u00d798=\"],\n  [\"51\u00d739=\", \"43\u00d769=\"],\n  [\"43\u00d757=\", \"18\u00d760=\"],\n  [\"40\u00d792=\", \"49\u00d778=\"],\n  [\"18\u00d794=\", \"99\u00d752=\"],\n  [\"94\u00d780=\", \"54\u00d757=\"],\n  [\"34\u00d793=\", \"83\u00d720=\"],\n  [\"84\u00d783=\", \"61\u00d736=\"],\n  [\"61\u00d737=\", \"55\u00d744=\"],\n  [\"66\u00d712=\", \"11\u00d752=\"],\n  [\"53\u00d712=\", \"46\u00d748=\"],\n  [\"76\u00d784=\", \"78\u00d714=\"],\n  [\"60\u00d722=\", \"95\u00d722=\"],\n  [\"18\u00d714=\", \"87\u00d740=\"],\n  [\"84\u00d739=\", \"99\u00d736=\"],\n  [\"38\u00d729=\", \"30\u00d784=\"],\n  [\"94\u00d730=\", \"33\u00d747=\"],\n  [\"32\u00d722=\", \"82\u00d717=\"],\n  [\"93\u00d751=\", \"16\u00d743=\"],\n  [\"32\u00d748=\", \"99\u00d771=\"],\n  [\"83\u00d786=\", \"77\u00d728=\"],\n  [\"82\u00d722=\", \"97\u00d720=\"],\n  [\"60\u00d715=\", \"97\u00d786=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and all the two-digit-by-two-digit multiplication\n# prompts in the table with the new values from the commit.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @{old=\"2025-05-14 Wednesday\"; new=\"2025-05-15 Thursday\"},\n    @{old=\"17\u00d780=\"; new=\"99\u00d799=\"},\n    @{old=\"35\u00d762=\"; new=\"84\u00d717=\"},\n    @{old=\"74\u00d772=\"; new=\"39\u00d798=\"},\n    @{old=\"51\u00d739=\"; new=\"43\u00d769=\"},\n    @{old=\"43\u00d757=\"; new=\"18\u00d760=\"},\n    @{old=\"40\u00d792=\"; new=\"49\u00d778=\"},\n    @{old=\"18\u00d794=\"; new=\"99\u00d752=\"},\n    @{old=\"94\u00d780=\"; new=\"54\u00d757=\"},\n    @{old=\"34\u00d793=\"; new=\"83\u00d720=\"},\n    @{old=\"84\u00d783=\"; new=\"61\u00d736=\"},\n    @{old=\"61\u00d737=\"; new=\"55\u00d744=\"},\n    @{old=\"66\u00d712=\"; new=\"11\u00d752=\"},\n    @{old=\"53\u00d712=\"; new=\"46\u00d748=\"},\n    @{old=\"76\u00d784=\"; new=\"78\u00d714=\"},\n    @{old=\"60\u00d722=\"; new=\"95\u00d722=\"},\n    @{old=\"18\u00d714=\"; new=\"87\u00d740=\"},\n    @{old=\"84\u00d739=\"; new=\"99\u00d736=\"},\n    @{old=\"38\u00d729=\"; new=\"30\u00d784=\"},\n    @{old=\"94\u00d730=\"; new=\"33\u00d747=\"},\n    @{old=\"32\u00d722=\"; new=\"82\u00d717=\"},\n    @{old=\"93\u00d751=\"; new=\"16\u00d743=\"},\n    @{old=\"32\u00d748=\"; new=\"99\u00d771=\"},\n    @{old=\"83\u00d786=\"; new=\"77\u00d728=\"},\n    @{old=\"82\u00d722=\"; new=\"97\u00d720=\"},\n    @{old=\"60\u00d715=\"; new=\"97\u00d786=\"}\n)\n\nforeach ($p in $pairs) {\n    $find = $d.Content.Find\n    $find.Text = $p.old\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Replacement.Text = $p.new\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
